$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Soda" scene stats (row 10) after raising MAX_VMAP_SIZE 64K -> 256K
$ws.Range("L10").Value = 525
$ws.Range("M10").Value = 254
$ws.Range("N10").Value = 102

# Move the active selection to L12
$ws.Range("L12").Select()
